$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in results/profit for rows whose matches have now finished ---
$ws.Range("G30").Value = "Fallo"
$ws.Range("H30").Value = -1

$ws.Range("G49").Value = "Fallo"
$ws.Range("H49").Value = -1

$ws.Range("G61").Value = "Acierto"
$ws.Range("H61").Value = 0.53

$ws.Range("G69").Value = "Acierto"
$ws.Range("H69").Value = 2

$ws.Range("G72").Value = "Fallo"
$ws.Range("H72").Value = -1

$ws.Range("G78").Value = "Fallo"
$ws.Range("H78").Value = -1

# --- Append two new pending picks (rows 79-80), cloning the blank-result
#     layout of an existing pending row so the new "resultado"/"profit"
#     cells stay genuinely empty (not yet settled) like the source rows ---
$ws.Range("A71:H71").Copy($ws.Range("A79:H80"))

$ws.Range("A79").Value = 14656375
$ws.Range("B79").Value = "'2025-09-10"
$ws.Range("B79").Style = "Normal"
$ws.Range("C79").Value = "Volodymyr Iakubenko"
$ws.Range("D79").Value = "Liam Broady"
$ws.Range("E79").Value = "Gana Volodymyr Iakubenko"
$ws.Range("F79").Value = 7.5

$ws.Range("A80").Value = 14656379
$ws.Range("B80").Value = "'2025-09-10"
$ws.Range("B80").Style = "Normal"
$ws.Range("C80").Value = "Kris van Wyk"
$ws.Range("D80").Value = "Karim Mabrouk"
$ws.Range("E80").Value = "Gana Karim Mabrouk"
$ws.Range("F80").Value = 6.5
